# export-structure-usagers.xlsx
# Add a new column "Numéro de distribution spéciale (BP, TSA, etc)" right
# after "Statut de la domiciliation" (i.e. insert a new column K), shifting
# every following column one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column at K (11th column). Excel shifts the existing
# column K ("Statut de la domiciliation") and everything after it one
# column to the right, carrying over formatting/column widths.
$ws.Columns.Item(11).Insert()

# Header cell for the freshly inserted column.
$headerCell = $ws.Cells.Item(2, 11)
$headerCell.Value = "Numéro de distribution spéciale (BP, TSA, etc)"

# Match the look of the other header cells (bold, vertically centered,
# wrapped text).
$headerCell.Font.Bold = $true
$headerCell.Font.Name = "Calibri"
$headerCell.Font.Size = 12
$headerCell.VerticalAlignment = -4108
$headerCell.WrapText = $true

# Give the new column roughly the same width as its neighbours.
$ws.Columns.Item(11).ColumnWidth = 21

# Restore the active selection on the frozen pane.
$ws.Range("K5").Select()
